$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General-looking numeric strings) for Price column cells so
# they keep their original formatted text instead of being parsed as numbers.
$ws.Range("D2:D5").NumberFormat = "@"
$ws.Range("D7:D15").NumberFormat = "@"
$ws.Range("D17:D26").NumberFormat = "@"
$ws.Range("D28:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.896.55'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '1.873.77'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '301.51'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = '0.5318'
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("D8").Value = '0.3755'
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").Value = '0.07166'
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").Value = '21.59'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = '0.8848'
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("D12").Value = '0.08100'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '1.868.87'
$ws.Range("E13").Value = '  -1.59%  '
$ws.Range("D14").Value = '93.16'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '5.275'
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '14.72'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '0.000008546'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '26.920.18'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").Value = '4.983'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").Value = '10.69'
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").Value = '6.375'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").Value = '147.53'
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").Value = '2.257'
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").Value = '1.730'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").Value = '114.41'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '4.744'
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("D30").Value = '4.572'
$ws.Range("E30").Value = '  -6.44%  '
$ws.Range("D31").Value = '0.09112'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '0.8000'
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").Value = '0.04980'
$ws.Range("E33").Value = '  -1.16%  '
$ws.Range("D34").Value = '1.175'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").Value = '2.988'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '2.633'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").Value = '3.176'
$ws.Range("E37").Value = '  -5.01%  '
$ws.Range("D38").Value = '3.176'
$ws.Range("E38").Value = '  -5.01%  '
$ws.Range("D39").Value = '0.01952'
$ws.Range("E39").Value = '  -1.93%  '
$ws.Range("D40").Value = '1.068'
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("D41").Value = '6.656'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").Value = '8.867'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("D43").Value = '115.64'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '0.5035'
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("D45").Value = '0.1494'
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("D46").Value = '1.0000'
$ws.Range("E46").Value = '  -0.27%  '
$ws.Range("D47").Value = '9.923'
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").Value = '1.613'
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").Value = '37.73'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("D50").Value = '0.06048'
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("D51").Value = '62.29'
$ws.Range("E51").Value = '  -2.74%  '
